$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2190.2273
$ws.Range("I70").Value = 1400.2
$ws.Range("J70").Value = 2848.5833
$ws.Range("K70").Value = 4200.6
$ws.Range("L70").Value = 8545.749899999999
$ws.Range("M70").Value = -3930.6
$ws.Range("N70").Value = -9085.749899999999

$ws.Range("H73").Value = 2190.2273
$ws.Range("I73").Value = 1400.2
$ws.Range("J73").Value = 2848.5833
$ws.Range("K73").Value = 4200.6
$ws.Range("L73").Value = 8545.749899999999
$ws.Range("M73").Value = -3264.6
$ws.Range("N73").Value = -10417.7499

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H108").Value = 40000
$ws.Range("J108").Value = 40000
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680

$ws.Range("H110").Value = 39000
$ws.Range("J110").Value = 39000
$ws.Range("L110").Value = 39000
$ws.Range("N110").Value = -47180

$ws.Range("H111").Value = 1849.5
$ws.Range("I111").Value = 1849.5
$ws.Range("K111").Value = 5548.5
$ws.Range("M111").Value = -2481.5

$ws.Range("H112").Value = 1332.0201
$ws.Range("J112").Value = 1345.0515
$ws.Range("L112").Value = 4035.1545
$ws.Range("N112").Value = -6251.154500000001

$ws.Range("H113").Value = 3567.6086
$ws.Range("I113").Value = 3569.7058
$ws.Range("J113").Value = 3561.6667
$ws.Range("K113").Value = 3569.7058
$ws.Range("L113").Value = 3561.6667
$ws.Range("M113").Value = -315.7058000000002
$ws.Range("N113").Value = -10069.6667

$ws.Range("H114").Value = 40000
$ws.Range("J114").Value = 40000
$ws.Range("L114").Value = 40000
$ws.Range("N114").Value = -48678

$ws.Range("H138").Value = 1886.53
$ws.Range("I138").Value = 609.525
$ws.Range("J138").Value = 2737.8667
$ws.Range("K138").Value = 1828.575
$ws.Range("L138").Value = 8213.6001
$ws.Range("M138").Value = 3311.425
$ws.Range("N138").Value = -18493.6001

$ws.Range("H141").Value = 360048.72
$ws.Range("I141").Value = 1411.3077
$ws.Range("K141").Value = 4233.9231
$ws.Range("M141").Value = 946.0769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1033.46
$ws.Range("I32").Value = 931.75
$ws.Range("K32").Value = 931.75
$ws.Range("M32").Value = -644.75

$ws.Range("H61").Value = 2000.2972
$ws.Range("I61").Value = 785.03845
$ws.Range("J61").Value = 4872.727
$ws.Range("K61").Value = 785.03845
$ws.Range("L61").Value = 4872.727
$ws.Range("M61").Value = -573.03845
$ws.Range("N61").Value = -5296.727

$ws.Range("H136").Value = 2000.2972
$ws.Range("I136").Value = 785.03845
$ws.Range("J136").Value = 4872.727
$ws.Range("K136").Value = 2355.11535
$ws.Range("L136").Value = 14618.181
$ws.Range("M136").Value = 194.88465
$ws.Range("N136").Value = -19718.181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1558.8182
$ws.Range("I105").Value = 1336.5
$ws.Range("J105").Value = 1900.8462
$ws.Range("K105").Value = 1336.5
$ws.Range("L105").Value = 1900.8462
$ws.Range("M105").Value = 410.5
$ws.Range("N105").Value = -5394.8462

$ws.Range("H107").Value = 4265.0713
$ws.Range("I107").Value = 3213.875
$ws.Range("J107").Value = 5666.6665
$ws.Range("K107").Value = 3213.875
$ws.Range("L107").Value = 5666.6665
$ws.Range("M107").Value = -1293.875
$ws.Range("N107").Value = -9506.666499999999

$ws.Range("H134").Value = 2442.5527
$ws.Range("I134").Value = 1823.3077
$ws.Range("J134").Value = 3784.25
$ws.Range("K134").Value = 5469.9231
$ws.Range("L134").Value = 11352.75
$ws.Range("M134").Value = -2934.9231
$ws.Range("N134").Value = -16422.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3963.9092
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3963.9092
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 11891.7276
$ws.Range("N80").Value = -13763.7276
$ws.Range("M80").ClearContents()

$ws.Range("H83").Value = 3963.9092
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3963.9092
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 35675.1828
$ws.Range("N83").Value = -45035.1828
$ws.Range("M83").ClearContents()

$ws.Range("H110").Value = 3008.25
$ws.Range("I110").Value = 1959.8
$ws.Range("J110").Value = 3757.1428
$ws.Range("K110").Value = 5879.4
$ws.Range("L110").Value = 11271.4284
$ws.Range("M110").Value = -1789.4
$ws.Range("N110").Value = -19451.4284

$ws.Range("H112").Value = 1465
$ws.Range("I112").Value = 740.7143
$ws.Range("K112").Value = 2222.1429
$ws.Range("M112").Value = -1114.1429

$ws.Range("H131").Value = 1400.5
$ws.Range("I131").Value = 2240.8333
$ws.Range("J131").Value = 1160.4048
$ws.Range("K131").Value = 6722.499899999999
$ws.Range("L131").Value = 3481.2144
$ws.Range("M131").Value = -1682.499899999999
$ws.Range("N131").Value = -13561.2144

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 9832.429
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 9832.429
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 9832.429
$ws.Range("N46").Value = -10144.429
$ws.Range("M46").ClearContents()

$ws.Range("H126").Value = 3126.1072
$ws.Range("I126").Value = 2041.1
$ws.Range("J126").Value = 3728.889
$ws.Range("K126").Value = 6123.299999999999
$ws.Range("L126").Value = 11186.667
$ws.Range("M126").Value = -3653.299999999999
$ws.Range("N126").Value = -16126.667

$ws.Range("H138").Value = 36380.625
$ws.Range("J138").Value = 36380.625
$ws.Range("L138").Value = 36380.625
$ws.Range("N138").Value = -46660.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 70014
$ws.Range("J43").Value = 70014
$ws.Range("L43").Value = 70014
$ws.Range("N43").Value = -70400

$ws.Range("H135").Value = 30013.182
$ws.Range("J135").Value = 30013.182
$ws.Range("L135").Value = 30013.182
$ws.Range("N135").Value = -40153.182

$ws.Range("H137").Value = 32300
$ws.Range("J137").Value = 32300
$ws.Range("L137").Value = 32300
$ws.Range("N137").Value = -42500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()

$ws.Range("H45").Value = 13287.75
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 13287.75
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 13287.75
$ws.Range("N45").Value = -14269.75
$ws.Range("M45").ClearContents()

$ws.Range("H126").Value = 4002172.2
$ws.Range("I126").Value = 1016.55554
$ws.Range("J126").Value = 14290858
$ws.Range("K126").Value = 3049.66662
$ws.Range("L126").Value = 42872574
$ws.Range("M126").Value = -579.66662
$ws.Range("N126").Value = -42877514
